$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2 = "2025-10-19T23:56:04.381727"
    3 = "2025-10-19T23:56:04.381727"
    4 = "2025-10-19T23:56:04.381727"
    5 = "2025-10-19T23:56:04.382724"
    6 = "2025-10-19T23:56:04.382724"
    7 = "2025-10-19T23:56:04.382724"
    8 = "2025-10-19T23:56:04.382724"
    9 = "2025-10-19T23:56:04.382724"
    10 = "2025-10-19T23:56:04.382724"
    11 = "2025-10-19T23:56:04.382724"
    12 = "2025-10-19T23:56:04.382724"
    13 = "2025-10-19T23:56:04.383723"
    14 = "2025-10-19T23:56:04.383723"
    15 = "2025-10-19T23:56:04.383723"
    16 = "2025-10-19T23:56:04.383723"
    17 = "2025-10-19T23:56:04.383723"
    18 = "2025-10-19T23:56:04.383723"
    19 = "2025-10-19T23:56:04.383723"
    20 = "2025-10-19T23:56:04.383723"
    21 = "2025-10-19T23:56:04.383723"
    22 = "2025-10-19T23:56:04.384724"
    23 = "2025-10-19T23:56:04.384724"
    24 = "2025-10-19T23:56:04.384724"
    25 = "2025-10-19T23:56:04.384724"
    26 = "2025-10-19T23:56:04.384724"
    27 = "2025-10-19T23:56:04.384724"
    28 = "2025-10-19T23:56:04.384724"
    29 = "2025-10-19T23:56:04.384724"
    30 = "2025-10-19T23:56:04.384724"
    31 = "2025-10-19T23:56:04.385723"
    32 = "2025-10-19T23:56:04.385723"
    33 = "2025-10-19T23:56:04.385723"
    34 = "2025-10-19T23:56:04.385723"
    35 = "2025-10-19T23:56:04.385723"
    36 = "2025-10-19T23:56:04.385723"
    37 = "2025-10-19T23:56:04.385723"
    38 = "2025-10-19T23:56:04.385723"
    39 = "2025-10-19T23:56:04.385723"
    40 = "2025-10-19T23:56:04.386724"
    41 = "2025-10-19T23:56:04.386724"
    42 = "2025-10-19T23:56:04.386724"
    43 = "2025-10-19T23:56:04.386724"
    44 = "2025-10-19T23:56:04.386724"
    45 = "2025-10-19T23:56:04.386724"
    46 = "2025-10-19T23:56:04.482920"
    47 = "2025-10-19T23:56:04.482920"
    48 = "2025-10-19T23:56:04.482920"
    49 = "2025-10-19T23:56:04.482920"
    50 = "2025-10-19T23:56:04.483918"
    51 = "2025-10-19T23:56:04.483918"
    52 = "2025-10-19T23:56:04.484917"
    53 = "2025-10-19T23:56:04.484917"
    54 = "2025-10-19T23:56:04.484917"
    55 = "2025-10-19T23:56:04.484917"
    56 = "2025-10-19T23:56:04.485918"
    57 = "2025-10-19T23:56:04.485918"
    58 = "2025-10-19T23:56:04.485918"
    59 = "2025-10-19T23:56:04.486921"
    60 = "2025-10-19T23:56:04.486921"
    61 = "2025-10-19T23:56:04.486921"
    62 = "2025-10-19T23:56:04.486921"
    63 = "2025-10-19T23:56:04.487920"
    64 = "2025-10-19T23:56:04.487920"
    65 = "2025-10-19T23:56:04.487920"
    66 = "2025-10-19T23:56:04.487920"
    67 = "2025-10-19T23:56:04.488920"
    68 = "2025-10-19T23:56:04.488920"
    69 = "2025-10-19T23:56:04.488920"
    70 = "2025-10-19T23:56:04.488920"
    71 = "2025-10-19T23:56:04.489914"
    72 = "2025-10-19T23:56:04.489914"
    73 = "2025-10-19T23:56:04.489914"
    74 = "2025-10-19T23:56:04.489914"
    75 = "2025-10-19T23:56:04.585882"
    76 = "2025-10-19T23:56:04.585882"
    77 = "2025-10-19T23:56:04.586882"
    78 = "2025-10-19T23:56:04.586882"
    79 = "2025-10-19T23:56:04.586882"
    80 = "2025-10-19T23:56:04.586882"
    81 = "2025-10-19T23:56:04.587885"
    82 = "2025-10-19T23:56:04.587885"
    83 = "2025-10-19T23:56:04.587885"
    84 = "2025-10-19T23:56:04.587885"
    85 = "2025-10-19T23:56:04.588883"
    86 = "2025-10-19T23:56:04.588883"
    87 = "2025-10-19T23:56:04.588883"
    88 = "2025-10-19T23:56:04.588883"
    89 = "2025-10-19T23:56:04.588883"
    90 = "2025-10-19T23:56:04.589880"
    91 = "2025-10-19T23:56:04.589880"
    92 = "2025-10-19T23:56:04.589880"
    93 = "2025-10-19T23:56:04.589880"
    94 = "2025-10-19T23:56:04.590884"
    95 = "2025-10-19T23:56:04.590884"
    96 = "2025-10-19T23:56:04.590884"
    97 = "2025-10-19T23:56:04.590884"
    98 = "2025-10-19T23:56:04.592362"
    99 = "2025-10-19T23:56:04.592362"
    100 = "2025-10-19T23:56:04.592883"
    101 = "2025-10-19T23:56:04.592883"
    102 = "2025-10-19T23:56:04.592883"
    103 = "2025-10-19T23:56:04.667215"
    104 = "2025-10-19T23:56:04.667215"
    105 = "2025-10-19T23:56:04.668219"
    106 = "2025-10-19T23:56:04.668219"
    107 = "2025-10-19T23:56:04.668219"
    108 = "2025-10-19T23:56:04.668219"
    109 = "2025-10-19T23:56:04.668219"
    110 = "2025-10-19T23:56:04.669214"
    111 = "2025-10-19T23:56:04.669214"
    112 = "2025-10-19T23:56:04.669214"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
